$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold free-form text (e.g. "60.404.96",
# "  -0.95%  ") that must stay text. Excel's COM layer auto-coerces numeric-
# looking strings (e.g. "510.71") into real numbers when assigned straight to
# .Value. Force the target ranges to a text number-format first so every
# assignment below is stored verbatim as text, then restore the cell style so
# we don't leave a stray NumberFormat on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.404.96"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.608.04"
$ws.Range("E3").Value = "  -3.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "510.71"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "154.57"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").Value = "2.618.77"
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "3.065.00"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "60.369.58"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "21.62"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "2.612.45"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "350.99"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "10.61"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").Value = "6.15"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "60.53"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "0.423"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "0.0₃0842"
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "19.45"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "151.29"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "5.77"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "0.885"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "0.848"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").Value = "36.29"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "294.79"
$ws.Range("E42").Value = "  -6.43%  "
$ws.Range("D43").Value = "0.626"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").Value = "0.100"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "0.0556"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").Value = "19.90"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").Value = "10.31"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "2.000.32"
$ws.Range("E51").Value = "  -3.68%  "

# Restore the original (default) cell style so no stray text-format override
# is left behind on these cells.
$ws.Range("D2:E51").Style = "Normal"
